# Actualización automática 2025-07-03 14:35:08
#
# Updates the figures for client "PEÑAHERRERA MOLINA JORGE OCTAVIO"
# (advisor ALMEIDA CUATIN JHONATHANN CARLOS) to reflect a new PORCELANATO
# sale of 128.3 registered in julio, propagated across the three report
# sheets (VENTAS POR GRUPO, VENTA MENSUAL, CUMPLIMIENTO MENSUAL).

$wb = $excel.ActiveWorkbook

# --- Sheet 1: VENTAS POR GRUPO --------------------------------------------
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")

# PORCELANATO sale for the client in row 22 (was 0)
$wsGrupo.Range("M22").Value = 128.3

# TOTAL row counter of non-zero clients for PORCELANATO (3 -> 4 of 30)
$wsGrupo.Range("M32").Value = "4 de 30"

# --- Sheet 2: VENTA MENSUAL -------------------------------------------------
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")

# julio sale for the same client/row (was 0)
$wsMensual.Range("F22").Value = 128.3

# TOTAL row for julio (was 118.18, now + 128.3)
$wsMensual.Range("F32").Value = 246.48

# --- Sheet 3: CUMPLIMIENTO MENSUAL -----------------------------------------
$wsCumplimiento = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# PORCELANATO group row (row 15): VENTA, POR CUMPLIR, CUMPLIMIENTO
$wsCumplimiento.Range("D15").Value = 236.2
$wsCumplimiento.Range("E15").Value = 23222.62
$wsCumplimiento.Range("F15").Value = 0.01006870763320576

# TOTAL row (row 18): VENTA, POR CUMPLIR, CUMPLIMIENTO
$wsCumplimiento.Range("D18").Value = 236.2
$wsCumplimiento.Range("E18").Value = 33698.51607548726
$wsCumplimiento.Range("F18").Value = 0.006960423640338604
